# Fruta / hortaliza, semanal
# Swap the contents of row 2 <-> row 4, and row 3 <-> row 6
# for the columns D, K, L, M, N, O, P, R, S (the rest of the
# columns are identical between the swapped rows, so no visible
# change happens there).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D","K","L","M","N","O","P","R","S")

function Swap-Rows($ws, $rowA, $rowB, $cols) {
    foreach ($col in $cols) {
        $rangeA = $ws.Range("$col$rowA")
        $rangeB = $ws.Range("$col$rowB")
        $valA = $rangeA.Value()
        $valB = $rangeB.Value()
        $rangeA.Value = $valB
        $rangeB.Value = $valA
    }
}

Swap-Rows $ws 2 4 $cols
Swap-Rows $ws 3 6 $cols
